# Update cryptos list figures (prices / volume%) per scheduled refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("D2").Value = "59.142.37"
$ws.Range("E2").Value = "  +1.13%  "
# Row 3
$ws.Range("D3").Value = "2.649.15"
$ws.Range("E3").Value = "  +1.32%  "
# Row 4
$ws.Range("E4").Value = "  -0.44%  "
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "519.70"
$ws.Range("E5").Value = "  +2.73%  "
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "147.20"
$ws.Range("E6").Value = "  +1.16%  "
# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.996"
$ws.Range("E7").Value = "  -0.17%  "
# Row 8
$ws.Range("E8").Value = "  +1.14%  "
# Row 9
$ws.Range("D9").Value = "2.662.13"
$ws.Range("E9").Value = "  +0.31%  "
# Row 10
$ws.Range("E10").Value = "  -2.30%  "
# Row 11
$ws.Range("E11").Value = "  +1.09%  "
# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.337"
$ws.Range("E12").Value = "  +0.48%  "
# Row 13
$ws.Range("E13").Value = "  +0.66%  "
# Row 14
$ws.Range("D14").Value = "3.107.12"
$ws.Range("E14").Value = "  +0.55%  "
# Row 15
$ws.Range("D15").Value = "59.139.91"
$ws.Range("E15").Value = "  +1.12%  "
# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "20.99"
$ws.Range("E16").Value = "  -0.30%  "
# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000138"
$ws.Range("E17").Value = "  +0.75%  "
# Row 18
$ws.Range("D18").Value = "2.653.83"
$ws.Range("E18").Value = "  -0.09%  "
# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "351.24"
$ws.Range("E19").Value = "  +2.82%  "
# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.50"
$ws.Range("E20").Value = "  -1.33%  "
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.40"
$ws.Range("E21").Value = "  -0.06%  "
# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.21"
$ws.Range("E22").Value = "  +1.31%  "
# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.997"
$ws.Range("E23").Value = "  -0.29%  "
# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "61.97"
$ws.Range("E24").Value = "  +2.36%  "
# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.418"
$ws.Range("E25").Value = "  -0.71%  "
# Row 26
$ws.Range("E26").Value = "  +2.74%  "
# Row 27
$ws.Range("E27").Value = "  -0.20%  "
# Row 28
$ws.Range("D28").Value = "0.0₃0809"
$ws.Range("E28").Value = "  -0.60%  "
# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.12"
$ws.Range("E29").Value = "  +1.36%  "
# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.997"
$ws.Range("E30").Value = "  -0.12%  "
# Row 31
$ws.Range("E31").Value = "  -3.15%  "
# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "18.98"
$ws.Range("E32").Value = "  +0.75%  "
# Row 33
$ws.Range("E33").Value = "  +1.61%  "
# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "149.64"
$ws.Range("E34").Value = "  +0.57%  "
# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.966"
$ws.Range("E35").Value = "  -6.28%  "
# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.05"
$ws.Range("E36").Value = "  +1.62%  "
# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.18"
$ws.Range("E37").Value = "  +3.62%  "
# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.856"
$ws.Range("E38").Value = "  +0.35%  "
# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "36.65"
$ws.Range("E39").Value = "  +0.85%  "
# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.44"
$ws.Range("E40").Value = "  +2.56%  "
# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.66"
$ws.Range("E41").Value = "  +0.03%  "
# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "282.12"
$ws.Range("E42").Value = "  +1.53%  "
# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.997"
$ws.Range("E43").Value = "  -0.19%  "
# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0989"
$ws.Range("E44").Value = "  +0.16%  "
# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "19.86"
$ws.Range("E45").Value = "  +1.96%  "
# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.607"
$ws.Range("E46").Value = "  -2.60%  "
# Row 47
$ws.Range("D47").Value = "2.105.15"
$ws.Range("E47").Value = "  +7.05%  "
# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0529"
$ws.Range("E48").Value = "  -1.45%  "
# Row 49
$ws.Range("E49").Value = "  +1.12%  "
# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "10.31"
$ws.Range("E50").Value = "  +0.72%  "
# Row 51
$ws.Range("B51").Value = "InjectiveProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "19.04"
$ws.Range("E51").Value = "  +4.12%  "
